{"js": "// Remove the \"ACS Spring National Conference ... (virtual) ... March 2024\"\n// poster-presentation list item (the paragraph sits right after the\n// \"DOE GSP PI's Meeting ... April 2024\" entry).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targetStart = \"ACS Spring National Conference\";\nconst targetEnd = \"March 2024\";\n\nlet found = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const text = p.text || \"\";\n  if (\n    text.indexOf(targetStart) !== -1 &&\n    text.indexOf(\"(virtual)\") !== -1 &&\n    text.trim().endsWith(targetEnd)\n  ) {\n    found = p;\n    break;\n  }\n}\n\nif (!found) {\n  throw new Error(\"Target paragraph not found\");\n}\n\nfound.delete();\nawait context.sync();\n", "ps1": "# Remove the \"ACS Spring National Conference ... (virtual) ... March 2024\"\n# poster-presentation list item (the paragraph sits right after the\n# \"DOE GSP PI's Meeting ... April 2024\" entry).\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*ACS Spring National Conference*\" -and $t -like \"*(virtual)*\" -and $t -like \"*March 2024*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.Delete()\n}\n"}
